$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values after repulling data / recalculating mean
$updates = @{
    2  = 0
    5  = -4
    6  = -5
    13 = -5
    20 = 0
    27 = 6
    28 = 2
    30 = 5
    36 = -1
    37 = -5
    41 = 3
    42 = -1
    43 = -1
    46 = -6
    47 = 1
    48 = -3
    49 = -4
    54 = 0
    58 = -3
    60 = 2
    63 = -5
    64 = -12
    68 = -2
    70 = -2
    71 = 3
    73 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
